# Auto-generated Excel COM-interop edit script
# Applies updated market-price figures (currentAveragePrice / Leve price / profit columns)
# per the scheduled-runner data refresh, one worksheet at a time.

$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 1398.6
$ws.Range("I41").Value = 997.8570999999999
$ws.Range("K41").Value = 997.8570999999999
$ws.Range("M41").Value = -557.8570999999999
$ws.Range("H100").Value = 6641.7144
$ws.Range("I100").Value = 5877.7
$ws.Range("J100").Value = 7336.273
$ws.Range("K100").Value = 5877.7
$ws.Range("L100").Value = 7336.273
$ws.Range("M100").Value = -5336.7
$ws.Range("N100").Value = -8418.273000000001
$ws.Range("H101").Value = 999
$ws.Range("I101").Value = 998
$ws.Range("J101").Value = 999.5
$ws.Range("K101").Value = 2994
$ws.Range("L101").Value = 2998.5
$ws.Range("M101").Value = -1372
$ws.Range("N101").Value = -6242.5
$ws.Range("H107").Value = 2238.2222
$ws.Range("I107").Value = 3200
$ws.Range("J107").Value = 1468.8
$ws.Range("K107").Value = 3200
$ws.Range("L107").Value = 1468.8
$ws.Range("M107").Value = -1280
$ws.Range("N107").Value = -5308.8
$ws.Range("H138").Value = 3858.4
$ws.Range("I138").Value = 2909.3333
$ws.Range("J138").Value = 4634.909
$ws.Range("K138").Value = 8727.999899999999
$ws.Range("L138").Value = 13904.727
$ws.Range("M138").Value = -3587.999899999999
$ws.Range("N138").Value = -24184.727

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 19587.25
$ws.Range("I2").Value = 25838.928
$ws.Range("K2").Value = 25838.928
$ws.Range("M2").Value = -25725.928
$ws.Range("H61").Value = 3868.625
$ws.Range("I61").Value = 2821.5
$ws.Range("J61").Value = 7010
$ws.Range("K61").Value = 2821.5
$ws.Range("L61").Value = 7010
$ws.Range("M61").Value = -2609.5
$ws.Range("N61").Value = -7434
$ws.Range("H74").Value = 2189.7778
$ws.Range("I74").Value = 2369.75
$ws.Range("J74").Value = 750
$ws.Range("K74").Value = 2369.75
$ws.Range("L74").Value = 750
$ws.Range("M74").Value = -1495.75
$ws.Range("N74").Value = -2498
$ws.Range("H77").Value = 2189.7778
$ws.Range("I77").Value = 2369.75
$ws.Range("J77").Value = 750
$ws.Range("K77").Value = 11848.75
$ws.Range("L77").Value = 3750
$ws.Range("M77").Value = -7480.75
$ws.Range("N77").Value = -12486
$ws.Range("H102").Value = 5550
$ws.Range("I102").Value = 10000
$ws.Range("K102").Value = 10000
$ws.Range("M102").Value = -8378
$ws.Range("H116").Value = 19587.25
$ws.Range("I116").Value = 25838.928
$ws.Range("K116").Value = 25838.928
$ws.Range("M116").Value = -23544.928
$ws.Range("H122").Value = 5509.3447
$ws.Range("I122").Value = 5269.357
$ws.Range("J122").Value = 5733.3335
$ws.Range("K122").Value = 15808.071
$ws.Range("L122").Value = 17200.0005
$ws.Range("M122").Value = -13358.071
$ws.Range("N122").Value = -22100.0005
$ws.Range("H124").Value = 45000
$ws.Range("J124").Value = 45000
$ws.Range("L124").Value = 45000
$ws.Range("N124").Value = -54820
$ws.Range("H132").Value = 1165.8125
$ws.Range("I132").Value = 1199.9
$ws.Range("J132").Value = 654.5
$ws.Range("K132").Value = 3599.7
$ws.Range("L132").Value = 1963.5
$ws.Range("M132").Value = -1069.7
$ws.Range("N132").Value = -7023.5
$ws.Range("H136").Value = 3868.625
$ws.Range("I136").Value = 2821.5
$ws.Range("J136").Value = 7010
$ws.Range("K136").Value = 8464.5
$ws.Range("L136").Value = 21030
$ws.Range("M136").Value = -5914.5
$ws.Range("N136").Value = -26130

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 19587.25
$ws.Range("I3").Value = 25838.928
$ws.Range("K3").Value = 25838.928
$ws.Range("M3").Value = -25724.928
$ws.Range("H105").Value = 6667.5
$ws.Range("I105").Value = 4998.5
$ws.Range("K105").Value = 4998.5
$ws.Range("M105").Value = -3251.5
$ws.Range("H107").Value = 8241.9
$ws.Range("I107").Value = 8983.799999999999
$ws.Range("K107").Value = 8983.799999999999
$ws.Range("M107").Value = -7063.799999999999
$ws.Range("H134").Value = 9098.541999999999
$ws.Range("I134").Value = 3053.8823
$ws.Range("J134").Value = 23778.428
$ws.Range("K134").Value = 9161.6469
$ws.Range("L134").Value = 71335.284
$ws.Range("M134").Value = -6626.6469
$ws.Range("N134").Value = -76405.284

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3148.8
$ws.Range("I99").Value = 3148.8
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 3148.8
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -1650.8
$ws.Range("N99").ClearContents()
$ws.Range("H126").Value = 3148.8
$ws.Range("I126").Value = 3148.8
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 9446.400000000001
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -6976.400000000001
$ws.Range("N126").ClearContents()
$ws.Range("H134").Value = 3684.5557
$ws.Range("I134").Value = 4020.6667
$ws.Range("J134").Value = 2004
$ws.Range("K134").Value = 12062.0001
$ws.Range("L134").Value = 6012
$ws.Range("M134").Value = -9527.000100000001
$ws.Range("N134").Value = -11082

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 36714730
$ws.Range("I4").Value = 48952684
$ws.Range("J4").Value = 862.7143
$ws.Range("K4").Value = 146858052
$ws.Range("L4").Value = 2588.1429
$ws.Range("M4").Value = -146857940
$ws.Range("N4").Value = -2812.1429
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("H68").Value = 1482.6842
$ws.Range("J68").Value = 1548.6875
$ws.Range("L68").Value = 4646.0625
$ws.Range("N68").Value = -6268.0625
$ws.Range("H71").Value = 1482.6842
$ws.Range("J71").Value = 1548.6875
$ws.Range("L71").Value = 13938.1875
$ws.Range("N71").Value = -22050.1875
$ws.Range("H88").Value = 4937.375
$ws.Range("J88").Value = 5749.75
$ws.Range("L88").Value = 17249.25
$ws.Range("N88").Value = -18105.25
$ws.Range("H91").Value = 4937.375
$ws.Range("J91").Value = 5749.75
$ws.Range("L91").Value = 17249.25
$ws.Range("N91").Value = -20213.25
$ws.Range("I100").Value = 4999
$ws.Range("J100").Value = 5000
$ws.Range("K100").Value = 14997
$ws.Range("L100").Value = 15000
$ws.Range("M100").Value = -14186
$ws.Range("N100").Value = -16622

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 1002.2857
$ws.Range("I6").Value = 801.75
$ws.Range("J6").Value = 1269.6666
$ws.Range("K6").Value = 801.75
$ws.Range("L6").Value = 1269.6666
$ws.Range("M6").Value = -688.75
$ws.Range("N6").Value = -1495.6666
$ws.Range("H16").Value = 1002.2857
$ws.Range("I16").Value = 801.75
$ws.Range("J16").Value = 1269.6666
$ws.Range("K16").Value = 801.75
$ws.Range("L16").Value = 1269.6666
$ws.Range("M16").Value = -551.75
$ws.Range("N16").Value = -1769.6666
$ws.Range("H23").Value = 1804.6666
$ws.Range("J23").Value = 1804.6666
$ws.Range("L23").Value = 1804.6666
$ws.Range("N23").Value = -2250.6666
$ws.Range("H80").Value = 63744
$ws.Range("I80").Value = 125157.555
$ws.Range("K80").Value = 125157.555
$ws.Range("M80").Value = -124159.555
$ws.Range("H83").Value = 63744
$ws.Range("I83").Value = 125157.555
$ws.Range("K83").Value = 625787.7749999999
$ws.Range("M83").Value = -620795.7749999999
$ws.Range("H97").Value = 7700.778
$ws.Range("I97").Value = 3295.2632
$ws.Range("K97").Value = 3295.2632
$ws.Range("M97").Value = -2799.2632
$ws.Range("H102").Value = 3923.4
$ws.Range("I102").Value = 3936.8
$ws.Range("J102").Value = 3869.8
$ws.Range("K102").Value = 3936.8
$ws.Range("L102").Value = 3869.8
$ws.Range("M102").Value = -2314.8
$ws.Range("N102").Value = -7113.8
$ws.Range("H113").Value = 12674.154
$ws.Range("I113").Value = 12979.444
$ws.Range("K113").Value = 12979.444
$ws.Range("M113").Value = -10809.444
$ws.Range("H122").Value = 6274.64
$ws.Range("I122").Value = 5957.8667
$ws.Range("K122").Value = 17873.6001
$ws.Range("M122").Value = -15423.6001
$ws.Range("H132").Value = 6913.5117
$ws.Range("I132").Value = 6569.6665
$ws.Range("J132").Value = 7707
$ws.Range("K132").Value = 19708.9995
$ws.Range("L132").Value = 23121
$ws.Range("M132").Value = -17178.9995
$ws.Range("N132").Value = -28181

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 1668.8182
$ws.Range("I9").Value = 1884.1111
$ws.Range("J9").Value = 700
$ws.Range("K9").Value = 1884.1111
$ws.Range("L9").Value = 700
$ws.Range("M9").Value = -1660.1111
$ws.Range("N9").Value = -1148
$ws.Range("H22").Value = 2385.9722
$ws.Range("I22").Value = 2451.8845
$ws.Range("K22").Value = 2451.8845
$ws.Range("M22").Value = -2156.8845
$ws.Range("H27").Value = 2385.9722
$ws.Range("I27").Value = 2451.8845
$ws.Range("K27").Value = 2451.8845
$ws.Range("M27").Value = -2344.8845
$ws.Range("H134").Value = 63999
$ws.Range("J134").Value = 63999
$ws.Range("L134").Value = 63999
$ws.Range("N134").Value = -74139
$ws.Range("H140").Value = 75999
$ws.Range("J140").Value = 75999
$ws.Range("L140").Value = 75999
$ws.Range("N140").Value = -86359

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1889.4736
$ws.Range("I126").Value = 1837.9333
$ws.Range("J126").Value = 2082.75
$ws.Range("K126").Value = 5513.7999
$ws.Range("L126").Value = 6248.25
$ws.Range("M126").Value = -3043.7999
$ws.Range("N126").Value = -11188.25
$ws.Range("H138").Value = 140000
$ws.Range("J138").Value = 140000
$ws.Range("L138").Value = 140000
$ws.Range("N138").Value = -150280
